$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix H10 style (pick up style "2" like the rest of the row, via format copy) ---
$null = $ws.Range("G10").Copy()
$null = $ws.Range("H10").PasteSpecial(-4122)

# --- Seed new shared strings in the exact order needed so sharedStrings.xml
#     gets indices 223..229 assigned as: GTV Media Group Inc., G-Coins,
#     "New York " (trailing space), RvT tokens, Rivetz Corp.,
#     Unregistered Offering, BitConnect ---
$ws.Range("H11").Value = "GTV Media Group Inc."
$ws.Range("G11").Value = "G-Coins"
$ws.Range("M11").Value = "New York "
$ws.Range("G12").Value = "RvT tokens"
$ws.Range("H12").Value = "Rivetz Corp."
$ws.Range("E12").Value = "Unregistered Offering"
$ws.Range("H13").Value = "BitConnect"

# --- Row 11: GTV Media Group Inc. / G-Coins ---
$ws.Range("E11").Value = "Fraud"
$ws.Range("F11").Value = "Civil"
$ws.Range("I11").Value = "N/A"
$ws.Range("J11").Value = 539000000
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0

# --- Row 12: Rivetz Corp. / RvT tokens ---
$ws.Range("F12").Value = "Civil"
$ws.Range("I12").Value = "Ethereum"
$ws.Range("J12").Value = 18000000
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0
$ws.Range("M12").Value = "New York"

# --- Row 13: BitConnect ---
$ws.Range("E13").Value = "Unregistered Offering"
$ws.Range("F13").Value = "Civil"
$ws.Range("G13").Value = "N/A"
$ws.Range("I13").Value = "Bitcoin"
$ws.Range("J13").Value = 2000000
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = "New York"

# --- Row 14: BitConnect (duplicate of row 13) ---
$ws.Range("E14").Value = "Unregistered Offering"
$ws.Range("F14").Value = "Civil"
$ws.Range("G14").Value = "N/A"
$ws.Range("H14").Value = "BitConnect"
$ws.Range("I14").Value = "Bitcoin"
$ws.Range("J14").Value = 2000000
$ws.Range("K14").Value = 1
$ws.Range("L14").Value = 1
$ws.Range("M14").Value = "New York"

# --- Column widths (closest achievable via the character-width COM property) ---
$ws.Columns.Item(5).ColumnWidth = 21.75
$ws.Columns.Item(6).ColumnWidth = 19.42
$ws.Columns.Item(13).ColumnWidth = 19.42

# --- Selection moved to L14 ---
$null = $ws.Range("L14").Select()

Write-Host "done"
